# StockDataRefreshMessage: synchronize product_stock_data row updates
# triggered from the stock refresh views.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("1P" - Lapicero Retractil kilometrico): product color corrected
# rojo -> azul, and the numeric quantity/price columns were round-tripped
# through the UI as text (inline string) values. Preserve the numbers,
# just change their stored type to text.
$ws.Range("C7").Value = "Lapicero Retractil kilometrico azul"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "45"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1200"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2500"

$ws.Range("G7").Value = 45808.68072812058

# Row 9 ("1AP" - Jabon piel dove en barra x 5): quantity/price columns were
# previously stored as text (inline string); refresh message normalizes
# them back to numeric values.
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 5600
$ws.Range("F9").Value = 8500

$ws.Range("G9").Value = 45808.67727196759
